$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 29 with the new protocol_028 test case data
$ws.Range("A29").Value = "protocol_028"
$ws.Range("B29").Value = "y"
$ws.Range("C29").Value = "设置全局变量后验证登录"
$ws.Range("D29").Value = "Protocol"
$ws.Range("H29").Value = "set global wait_timeout=60"
$ws.Range("K29").Value = "connection"

# Update the K-column data validation dropdown list to include the new "connection" option
$dv = $ws.Range("K2:K1048576").Validation
$dv.Modify(3, 1, 1, '"csv_equals,csv_containsAll,string_equals,effected_rows_assert,table_assert,assertNull,justExec,connection,SQLException"')

# Move the active selection to match the post-edit cursor position
$ws.Range("J31").Select()
